# Insert a new data row at row 481 (pushing the existing rows 481-518 down
# to 482-519) and populate it with the new "Femacal de La Calera" / "Acelga"
# observation dated 45013 (2023-03-28).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 481..518 down to 482..519, leaving a blank row 481.
$ws.Rows.Item(481).Insert()

# Fill the new row with its data.
$ws.Cells.Item(481, 1).Value = 3
$ws.Cells.Item(481, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(481, 3).Value = "Coquimbo"
$ws.Cells.Item(481, 4).Value = 45013
$ws.Cells.Item(481, 5).Value = 5
$ws.Cells.Item(481, 6).Value = 100112009
$ws.Cells.Item(481, 7).Value = "Acelga"
$ws.Cells.Item(481, 8).Value = "Sin especificar"
$ws.Cells.Item(481, 9).Value = "Primera"
$ws.Cells.Item(481, 10).Value = 200
$ws.Cells.Item(481, 11).Value = 3500
$ws.Cells.Item(481, 12).Value = 3700
$ws.Cells.Item(481, 13).Value = 3590
$ws.Cells.Item(481, 14).Value = "$/docena de atados (6 kilos)"
$ws.Cells.Item(481, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(481, 16).Value = 598
$ws.Cells.Item(481, 17).Value = 6
$ws.Cells.Item(481, 18).Value = "Hortaliza"

Write-Output "inserted row 481"
